$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.522813
$ws.Range("H2").Value = 13.568439
$ws.Range("I2").Value = 0.5686345655850849
$ws.Range("J2").Value = 0.5686345655850847
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.522813
$ws.Range("N2").Value = 13.568439
$ws.Range("O2").Value = 0.5686345655850849
$ws.Range("P2").Value = 0.5686345655850847
$ws.Range("Q2").Value = 20.455837432969
$ws.Range("R2").Value = 184.102536896721
$ws.Range("S2").Value = 0.3233452691781382
$ws.Range("T2").Value = 0.3233452691781381

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.522813
$ws.Range("H3").Value = 13.568439
$ws.Range("I3").Value = 0.5686345655850849
$ws.Range("J3").Value = 0.5686345655850847
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.319478666666667
$ws.Range("N3").Value = 3.958436
$ws.Range("O3").Value = 0.1658925934852462
$ws.Range("P3").Value = 0.1658925934852462
$ws.Range("Q3").Value = 5.967755266822667
$ws.Range("R3").Value = 53.709797401404
$ws.Range("S3").Value = 0.09433226283026605
$ws.Range("T3").Value = 0.09433226283026605

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.522813
$ws.Range("H4").Value = 13.568439
$ws.Range("I4").Value = 0.5686345655850849
$ws.Range("J4").Value = 0.5686345655850847
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.111521333333334
$ws.Range("N4").Value = 6.334564
$ws.Range("O4").Value = 0.265472840929669
$ws.Range("P4").Value = 0.265472840929669
$ws.Range("Q4").Value = 9.550016136177335
$ws.Range("R4").Value = 85.950145225596
$ws.Range("S4").Value = 0.1509570335766807
$ws.Range("T4").Value = 0.1509570335766806

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.319478666666667
$ws.Range("H5").Value = 3.958436
$ws.Range("I5").Value = 0.1658925934852462
$ws.Range("J5").Value = 0.1658925934852462
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.522813
$ws.Range("N5").Value = 13.568439
$ws.Range("O5").Value = 0.5686345655850849
$ws.Range("P5").Value = 0.5686345655850847
$ws.Range("Q5").Value = 5.967755266822667
$ws.Range("R5").Value = 53.709797401404
$ws.Range("S5").Value = 0.09433226283026605
$ws.Range("T5").Value = 0.09433226283026605

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.319478666666667
$ws.Range("H6").Value = 3.958436
$ws.Range("I6").Value = 0.1658925934852462
$ws.Range("J6").Value = 0.1658925934852462
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.319478666666667
$ws.Range("N6").Value = 3.958436
$ws.Range("O6").Value = 0.1658925934852462
$ws.Range("P6").Value = 0.1658925934852462
$ws.Range("Q6").Value = 1.741023951788444
$ws.Range("R6").Value = 15.669215566096
$ws.Range("S6").Value = 0.02752035257326115
$ws.Range("T6").Value = 0.02752035257326116

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.319478666666667
$ws.Range("H7").Value = 3.958436
$ws.Range("I7").Value = 0.1658925934852462
$ws.Range("J7").Value = 0.1658925934852462
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.111521333333334
$ws.Range("N7").Value = 6.334564
$ws.Range("O7").Value = 0.265472840929669
$ws.Range("P7").Value = 0.265472840929669
$ws.Range("Q7").Value = 2.786107353544889
$ws.Range("R7").Value = 25.074966181904
$ws.Range("S7").Value = 0.04403997808171901
$ws.Range("T7").Value = 0.04403997808171902

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.111521333333334
$ws.Range("H8").Value = 6.334564
$ws.Range("I8").Value = 0.265472840929669
$ws.Range("J8").Value = 0.265472840929669
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.522813
$ws.Range("N8").Value = 13.568439
$ws.Range("O8").Value = 0.5686345655850849
$ws.Range("P8").Value = 0.5686345655850847
$ws.Range("Q8").Value = 9.550016136177335
$ws.Range("R8").Value = 85.950145225596
$ws.Range("S8").Value = 0.1509570335766807
$ws.Range("T8").Value = 0.1509570335766806

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.111521333333334
$ws.Range("H9").Value = 6.334564
$ws.Range("I9").Value = 0.265472840929669
$ws.Range("J9").Value = 0.265472840929669
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.319478666666667
$ws.Range("N9").Value = 3.958436
$ws.Range("O9").Value = 0.1658925934852462
$ws.Range("P9").Value = 0.1658925934852462
$ws.Range("Q9").Value = 2.786107353544889
$ws.Range("R9").Value = 25.074966181904
$ws.Range("S9").Value = 0.04403997808171901
$ws.Range("T9").Value = 0.04403997808171902

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.111521333333334
$ws.Range("H10").Value = 6.334564
$ws.Range("I10").Value = 0.265472840929669
$ws.Range("J10").Value = 0.265472840929669
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.111521333333334
$ws.Range("N10").Value = 6.334564
$ws.Range("O10").Value = 0.265472840929669
$ws.Range("P10").Value = 0.265472840929669
$ws.Range("Q10").Value = 4.458522341121779
$ws.Range("R10").Value = 40.126701070096
$ws.Range("S10").Value = 0.07047582927126934
$ws.Range("T10").Value = 0.07047582927126934
